# New weekly price entry for "Hortaliza, Macroferia Regional de Talca - Papa"
# A new record (Rodeo / 1a (guarda lavada), Region de La Araucania) is added
# as the new first data row of the "recent" block, row 248, pushing the
# existing rows 248-285 down to 249-286.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 248, shifting rows 248:285 down to 249:286
$ws.Rows.Item(248).Insert()

# Populate the new row 248 with the new record
$ws.Range("A248").Value = 5
$ws.Range("B248").Value = "Macroferia Regional de Talca"
$ws.Range("C248").Value = "Maule"
$ws.Range("D248").Value = 44474
$ws.Range("E248").Value = 7
$ws.Range("F248").Value = 100114001
$ws.Range("G248").Value = "Papa"
$ws.Range("H248").Value = "Rodeo"
$ws.Range("I248").Value = "1a (guarda lavada)"
$ws.Range("J248").Value = 1500
$ws.Range("K248").Value = 9000
$ws.Range("L248").Value = 9000
$ws.Range("M248").Value = 9000
$ws.Range("N248").Value = "`$/malla 25 kilos"
$ws.Range("O248").Value = "Región de La Araucanía"
$ws.Range("P248").Value = 360
$ws.Range("Q248").Value = 25
$ws.Range("R248").Value = "Hortaliza"

# Note: Rows.Insert() already carries the row-248 (now row-249) formatting
# up into the new row, including the date-style (s="2") on column D, so no
# extra style assignment is required here.
